$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2-398). The diff updates every one of these cells from the
# serial value 45190 to 45192, leaving everything else untouched.
$ws.Range("C2:C398").Value = 45192
